$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "iQOO Z10 5G (Stellar Black, 12GB RAM, 256GB Storage) | India's Biggest Ever 7300 mAh Battery | Snapdragon 7s Gen 3 Processor | Brightest Quad Curved AMOLED Display in The Segment"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "25,998"

$ws.Range("A8").Value = "realme NARZO 80 Pro 5G (Speed Silver,12GB+256GB) | Segment's 1st MediaTek Dimensity 7400 Chipset | 6000mAh Titan Battery + 80W Ultra Charge | 4500nits HyperGlow Esports Display | IP69 Waterproof"
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "21,498"

$ws.Range("A9").Value = "POCO C71, Desert Gold (6GB, 128GB)"
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "6,799"

$ws.Range("A10").Value = "Redmi 15 5G Midnight Black 8GB + 256GB | Segment's Largest 7000mAhA Battery | Segment's Largest Display 17.53cm(6.9) Up to 144Hz | Snapdragon 6s Gen 3 | 18W Reverse Charging | 50MP AI Dual Camera"
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "16,998"

$ws.Range("A11").Value = "Samsung Galaxy M06 5G (Sage Green, 6GB RAM, 128 GB Storage) | MediaTek Dimensity 6300 | AnTuTu Score 422K+ | 12 5G Bands| 25W Fast Charging | 4 Gen of OS Upgrades | Without Charger"
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "8,999"

$ws.Range("A12").Value = "Samsung Galaxy A55 5G (Awesome Iceblue, 8GB RAM, 128GB Storage) | AI | Metal Frame | 50 MP Main Camera (OIS) | Super HDR Video| Nightography | IP67 | Corning Gorilla Glass Victus+ | sAMOLED Display"
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "23,999"

$ws.Range("A21").Value = "Nokia All-New 105 Single Sim Keypad Phone with Built-in UPI Payments, Long-Lasting Battery, Wireless FM Radio | Red"
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = "1,199"
